$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.107.50'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '2.315.29'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.18'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.77%  '
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.117'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.03'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").Value = '2.676.11'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '2.288.35'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").Value = '43.014.63'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.54'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.90%  '
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.17'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.03'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.62'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.04'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.62'
$ws.Range("D28").ClearFormats()
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  -10.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.57'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.25'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.66%  '
$ws.Range("E33").Value = '  +6.23%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.36'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.80%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0695'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.102'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("E40").Value = '  -1.99%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").Value = '2.000.68'
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.17'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.83%  '
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.53'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.89'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.13'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +7.07%  '
$ws.Range("D50").Value = '2.541.78'
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("E51").Value = '  +0.33%  '
